# Apply the updates described by the commit diff to the workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: Metadata ---
$ws1 = $wb.Worksheets.Item(1)

# Version: 5.0.0 -> 6.0.0
$ws1.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank, now "Alvearie Team"
$ws1.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> "Jurisdiction" / "United States of America"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row -> remove it entirely
$ws1.Rows.Item(11).Delete()

# --- Sheet 2: Elements ---
$ws2 = $wb.Worksheets.Item(2)

# Root "Extension" element row (row 2): Short/Definition updated to match this profile
$ws2.Range("K2").Value = "Communication Outcome"
$ws2.Range("L2").Value = "Vendor-specific communication response code or text"
